# DPLKKPS140-001 - update Register number from M03220800000027 to M03220800000018
# and move the active selection to F2 (scrolled so column D is leftmost).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRegister = "M03220800000018"

# NO_REGISTER column (N2) holds just the register number.
$ws.Range("N2").Value = $newRegister

# PREPARATION column (F2) holds a multi-line note that also embeds the register number.
$nl = [char]10
$ws.Range("F2").Value = "Username : 31816;" + $nl + "Password : bni1234;" + $nl + "Role : 09 - Penyelia Settlement;" + $nl + "No Register : " + $newRegister

# Update the view/selection: scroll so column D is the left-most visible column
# and make F2 the active/selected cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F2").Select()
